$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A..T (Sending cluster, Ligand
# symbol, Receptor symbol, Target cluster, Ligand-expressing cells, Ligand
# detection rate, Ligand average/total expression value, Ligand derived
# specificity of average/total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average/total expression value, Receptor
# derived specificity of average/total expression value, Edge average/total
# expression weight, Edge average/total expression derived specificity).
$rows = @(
    @(2, "ECs", "Fn1", "Tshr", "ECs", 3, 1, 40.75339133333333, 122.260174, 0.02126536631186857, 0.02126536631186857, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 29.38258385046334, 264.4432546541701, 0.003260672918560337, 0.003260672918560337),
    @(3, "ECs", "Fn1", "Tshr", "FAPs", 3, 1, 40.75339133333333, 122.260174, 0.02126536631186857, 0.02126536631186857, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 138.0441118925015, 1242.397007032514, 0.01531916660241216, 0.01531916660241216),
    @(4, "ECs", "Fn1", "Tshr", "M2", 3, 1, 40.75339133333333, 122.260174, 0.02126536631186857, 0.02126536631186857, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 4.75794485370289, 42.82150368332601, 0.0005280033237181724, 0.0005280033237181723),
    @(5, "ECs", "Fn1", "Tshr", "sCs", 3, 1, 40.75339133333333, 122.260174, 0.02126536631186857, 0.02126536631186857, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 19.44188079179889, 174.97692712619, 0.002157523467177901, 0.0021575234671779),
    @(6, "FAPs", "Fn1", "Tshr", "ECs", 3, 1, 1689.289306666667, 5067.86792, 0.8814813868902838, 0.8814813868902838, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 1217.952250767067, 10961.5702569036, 0.1351597919497866, 0.1351597919497865),
    @(7, "FAPs", "Fn1", "Tshr", "FAPs", 3, 1, 1689.289306666667, 5067.86792, 0.8814813868902838, 0.8814813868902838, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 5722.13586253279, 51499.22276279511, 0.635002474194908, 0.635002474194908),
    @(8, "FAPs", "Fn1", "Tshr", "M2", 3, 1, 1689.289306666667, 5067.86792, 0.8814813868902838, 0.8814813868902838, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 197.2239634568978, 1775.01567111208, 0.02188653114402324, 0.02188653114402324),
    @(9, "FAPs", "Fn1", "Tshr", "sCs", 3, 1, 1689.289306666667, 5067.86792, 0.8814813868902838, 0.8814813868902838, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 805.8951721205779, 7253.0565490852, 0.08943258960156604, 0.08943258960156603),
    @(10, "M2", "Fn1", "Tshr", "ECs", 3, 1, 100.9654023333333, 302.896207, 0.05268435816499466, 0.05268435816499466, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 72.79454060129834, 655.1508654116851, 0.008078227169049719, 0.008078227169049717),
    @(11, "M2", "Fn1", "Tshr", "FAPs", 3, 1, 100.9654023333333, 302.896207, 0.05268435816499466, 0.05268435816499466, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 342.0004775301752, 3078.004297771577, 0.03795281248554185, 0.03795281248554185),
    @(12, "M2", "Fn1", "Tshr", "M2", 3, 1, 100.9654023333333, 302.896207, 0.05268435816499466, 0.05268435816499466, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 11.78767706728256, 106.089093605543, 0.001308113662897515, 0.001308113662897515),
    @(13, "M2", "Fn1", "Tshr", "sCs", 3, 1, 100.9654023333333, 302.896207, 0.05268435816499466, 0.05268435816499466, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 48.16672311281056, 433.5005080152951, 0.005345204847505575, 0.005345204847505575),
    @(14, "sCs", "Fn1", "Tshr", "ECs", 3, 1, 85.41274733333334, 256.238242, 0.04456888863285297, 0.04456888863285297, 3, 1, 0.7209850000000001, 2.162955, 0.1533325535399077, 0.1533325535399077, 61.58130963612335, 554.2317867251101, 0.006833861502511113, 0.006833861502511112),
    @(15, "sCs", "Fn1", "Tshr", "FAPs", 3, 1, 85.41274733333334, 256.238242, 0.04456888863285297, 0.04456888863285297, 3, 1, 3.387303666666666, 10.161911, 0.7203810354238886, 0.7203810354238885, 289.3189122222736, 2603.870210000462, 0.0321065821410266, 0.03210658214102659),
    @(16, "sCs", "Fn1", "Tshr", "M2", 3, 1, 85.41274733333334, 256.238242, 0.04456888863285297, 0.04456888863285297, 2, 0.6666666666666666, 0.1167496666666667, 0.350249, 0.02482926068494219, 0.02482926068494219, 9.971909780250892, 89.74718802225802, 0.001106612554303263, 0.001106612554303263),
    @(17, "sCs", "Fn1", "Tshr", "sCs", 3, 1, 85.41274733333334, 256.238242, 0.04456888863285297, 0.04456888863285297, 3, 1, 0.4770616666666667, 1.431185, 0.1014571503512615, 0.1014571503512615, 40.74714759741889, 366.7243283767701, 0.004521832435011993, 0.004521832435011993)
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}
